$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 6).Value = 361
$ws.Cells.Item(7, 6).Value = 1177
$ws.Cells.Item(9, 6).Value = 7116
$ws.Cells.Item(11, 6).Value = 89
$ws.Cells.Item(12, 6).Value = 2045
$ws.Cells.Item(13, 6).Value = 7986
$ws.Cells.Item(14, 6).Value = 44
$ws.Cells.Item(16, 6).Value = 5509
$ws.Cells.Item(17, 6).Value = 50
$ws.Cells.Item(18, 6).Value = 2414
$ws.Cells.Item(19, 6).Value = 1025
$ws.Cells.Item(21, 6).Value = 301
$ws.Cells.Item(26, 6).Value = 257
$ws.Cells.Item(28, 6).Value = 2362
$ws.Cells.Item(30, 6).Value = 264
$ws.Cells.Item(32, 6).Value = 147
$ws.Cells.Item(33, 6).Value = 584
$ws.Cells.Item(36, 6).Value = 1498
$ws.Cells.Item(37, 6).Value = 34
$ws.Cells.Item(39, 6).Value = 2322
$ws.Cells.Item(40, 6).Value = 2215
$ws.Cells.Item(42, 6).Value = 9

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 77
$ws.Cells.Item(5, 6).Value = 6

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 255
$ws.Cells.Item(3, 6).Value = 1279

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 255
$ws.Cells.Item(4, 6).Value = 1279
$ws.Cells.Item(8, 6).Value = 361
$ws.Cells.Item(9, 6).Value = 1177
$ws.Cells.Item(11, 6).Value = 7116
$ws.Cells.Item(13, 6).Value = 89
$ws.Cells.Item(14, 6).Value = 2045
$ws.Cells.Item(15, 6).Value = 7986
$ws.Cells.Item(16, 6).Value = 44
$ws.Cells.Item(18, 6).Value = 5509
$ws.Cells.Item(19, 6).Value = 50
$ws.Cells.Item(20, 6).Value = 2414
$ws.Cells.Item(21, 6).Value = 1025
$ws.Cells.Item(25, 6).Value = 77
$ws.Cells.Item(30, 6).Value = 2362
$ws.Cells.Item(32, 6).Value = 264
$ws.Cells.Item(34, 6).Value = 147
$ws.Cells.Item(35, 6).Value = 6
$ws.Cells.Item(36, 6).Value = 584
$ws.Cells.Item(40, 6).Value = 1498
$ws.Cells.Item(41, 6).Value = 34
$ws.Cells.Item(43, 6).Value = 2322
$ws.Cells.Item(45, 6).Value = 2215
$ws.Cells.Item(47, 6).Value = 9
